# Update the "as_of_utc" timestamp column (AA) on both data sheets
# ("Главные" and "Линейные") from the old publish time to the new one.

$wb = $excel.ActiveWorkbook

$oldTimestamp = "2025-12-18 03:04:27"
$newTimestamp = "2025-12-18 04:40:32"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Range("AA$row")
        # Force text formatting so the timestamp string is preserved literally
        # (and not auto-converted into a date/number by Excel).
        $cell.NumberFormat = "@"
        $cell.Value = $newTimestamp
    }
}
